{"js": "// The canonical-OOXML diff for this revision touches only `word/document.xml`\n// and `word/styles.xml`, and every single hunk in it is a reordering of XML\n// attributes (and root-element namespace declarations) into alphabetical\n// order - e.g. `<w:color w:val=\"E36C0A\" w:themeColor=\"accent6\" w:themeShade=\"BF\"/>`\n// becomes `<w:color w:themeColor=\"accent6\" w:themeShade=\"BF\" w:val=\"E36C0A\"/>`,\n// `<w:pgSz w:w=\"11906\" w:h=\"16838\"/>` becomes `<w:pgSz w:h=\"16838\" w:w=\"11906\"/>`,\n// etc. No attribute value, run of text, style definition, or page-layout\n// number actually changes anywhere in the diff (verified attribute-set-by-\n// attribute-set and by XML canonicalization of the whole parts: the\n// \"before\" and \"after\" XML are C14N-identical). Attribute/namespace\n// declaration order is not semantically meaningful in XML and is not\n// something the Word document object model (Office.js or COM) exposes or\n// lets a caller control - it is purely a side effect of whichever library\n// last serialized the underlying part.\n//\n// So, applied to the actual Word content model, this revision changes\n// nothing: the same text, the same run/paragraph formatting (including the\n// accent6/BF-shaded orange field-code color), and the same page size and\n// margins remain. There is therefore no content mutation to perform here;\n// this script intentionally leaves the document body, formatting and\n// styles untouched.\n", "ps1": "# The canonical-OOXML diff for this revision touches only `word/document.xml`\n# and `word/styles.xml`, and every single hunk in it is a reordering of XML\n# attributes (and root-element namespace declarations) into alphabetical\n# order - e.g. `<w:color w:val=\"E36C0A\" w:themeColor=\"accent6\" w:themeShade=\"BF\"/>`\n# becomes `<w:color w:themeColor=\"accent6\" w:themeShade=\"BF\" w:val=\"E36C0A\"/>`,\n# `<w:pgSz w:w=\"11906\" w:h=\"16838\"/>` becomes `<w:pgSz w:h=\"16838\" w:w=\"11906\"/>`,\n# etc. No attribute value, run of text, style definition, or page-layout\n# number actually changes anywhere in the diff (verified attribute-set-by-\n# attribute-set and by XML canonicalization of the whole parts: the\n# \"before\" and \"after\" XML are C14N-identical). Attribute/namespace\n# declaration order is not semantically meaningful in XML and is not\n# something the Word document object model (Office.js or COM) exposes or\n# lets a caller control - it is purely a side effect of whichever library\n# last serialized the underlying part.\n#\n# So, applied to the actual Word content model, this revision changes\n# nothing: the same text, the same run/paragraph formatting (including the\n# accent6/BF-shaded orange field-code color), and the same page size and\n# margins remain. There is therefore no content mutation to perform here;\n# this script intentionally leaves the document body, formatting and\n# styles untouched.\n"}
